# Apply the updated cryptocurrency price/volume figures.
# The "Price" (column D) and "Volume(1h)" (column E) cells hold
# plain text values (e.g. thousands-grouped prices like "27.310.44"
# and padded percentages like "  -2.00%  "). Several of the new
# price strings look like ordinary decimal numbers (e.g. "1.006"),
# so we force those cells to Text format first -- otherwise Excel
# would silently reinterpret them as numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.310.44"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "1.826.28"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -1.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.54"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3691"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07255"
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8647"
$ws.Range("E10").Value = "  -2.43%  "
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("D12").Value = "1.818.58"
$ws.Range("E12").Value = "  -2.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.711"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07103"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("E15").Value = "  -3.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.64"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008869"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("D21").Value = "27.325.99"
$ws.Range("E21").Value = "  -2.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.147"
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("D24").Value = "2.051.54"
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.13"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.35"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.148"
$ws.Range("E28").Value = "  +6.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.255"
$ws.Range("E29").Value = "  -3.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.61"
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08905"
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.202"
$ws.Range("E32").Value = "  -3.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7572"
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.460"
$ws.Range("E34").Value = "  -3.01%  "
$ws.Range("E35").Value = "  -2.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.005"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.113"
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01982"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05277"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.169"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.875"
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1698"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5052"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.678"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "107.74"
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4757"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.666"
$ws.Range("E50").Value = "  -3.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.833"
$ws.Range("E51").Value = "  -2.73%  "
